# Commit: "Update: Them tien mua USB UART FT232"
# Adds a new expense row (row 9) to the cost-tracking sheet for the
# purchase of a "USB Serial FT232" (USB-UART FT232RL converter board).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (row 8) down to the
# new row 9 so borders/alignment/number-formats stay consistent with the
# rest of the table.
$ws.Range("A8:D8").Copy()
$ws.Range("A9:D9").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the new expense entry.
$ws.Range("A9").Value = 43195   # Ngay: 2018-04-05 (Excel date serial)
$ws.Range("B9").Value = "USB Serial FT232"
$ws.Range("C9").Value = 60000
$ws.Range("D9").Value = "https://icdayroi.com/mach-chuyen-usb-uart-ttl-ft232rl"

# Move the selection to A10, matching where the cursor ends up after
# entering the new row.
$ws.Range("A10").Select()
